# Auto-generated edit script: numeric snapshot refresh across ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets
# (computed "profit" columns H-N recalculated; some cells added/removed as values became N/A or newly populated)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (45 cell(s)) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 624.36365
$ws.Range("I33").Value = 315.5
$ws.Range("K33").Value = 315.5
$ws.Range("M33").Value = -86.5
$ws.Range("H40").Value = 5434.364
$ws.Range("I40").Value = 6064.4443
$ws.Range("J40").Value = 2599
$ws.Range("K40").Value = 6064.4443
$ws.Range("L40").Value = 2599
$ws.Range("M40").Value = -5889.4443
$ws.Range("N40").Value = -2949
$ws.Range("H69").Value = 8831.75
$ws.Range("I69").Value = 8831.75
$ws.Range("K69").Value = 26495.25
$ws.Range("M69").Value = -25621.25
$ws.Range("H70").Value = 5431.857
$ws.Range("I70").Value = 4414.75
$ws.Range("J70").Value = 6788
$ws.Range("K70").Value = 13244.25
$ws.Range("L70").Value = 20364
$ws.Range("M70").Value = -12974.25
$ws.Range("N70").Value = -20904
$ws.Range("H72").Value = 8831.75
$ws.Range("I72").Value = 8831.75
$ws.Range("K72").Value = 79485.75
$ws.Range("M72").Value = -75117.75
$ws.Range("H73").Value = 5431.857
$ws.Range("I73").Value = 4414.75
$ws.Range("J73").Value = 6788
$ws.Range("K73").Value = 13244.25
$ws.Range("L73").Value = 20364
$ws.Range("M73").Value = -12308.25
$ws.Range("N73").Value = -22236
$ws.Range("H131").Value = 45051.375
$ws.Range("I131").Value = 1381.1333
$ws.Range("K131").Value = 4143.3999
$ws.Range("M131").Value = 896.6000999999997
$ws.Range("H137").Value = 38565.645
$ws.Range("I137").Value = 41036.32
$ws.Range("K137").Value = 123108.96
$ws.Range("M137").Value = -120558.96
$ws.Range("H138").Value = 3165.6326
$ws.Range("J138").Value = 4053.4768
$ws.Range("L138").Value = 12160.4304
$ws.Range("N138").Value = -22440.4304

# ---- Sheet: ARM (4 cell(s)) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5627.5713
$ws.Range("I122").Value = 2717.5454
$ws.Range("K122").Value = 8152.6362
$ws.Range("M122").Value = -5702.6362

# ---- Sheet: CRP (107 cell(s)) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 166672580
$ws.Range("I4").Value = 4999.5
$ws.Range("K4").Value = 4999.5
$ws.Range("M4").Value = -4887.5
$ws.Range("H16").Value = 1767.5
$ws.Range("I16").Value = 1981.3334
$ws.Range("J16").Value = 1639.2
$ws.Range("K16").Value = 1981.3334
$ws.Range("L16").Value = 1639.2
$ws.Range("M16").Value = -1694.3334
$ws.Range("N16").Value = -2213.2
$ws.Range("H19").Value = 1273.2858
$ws.Range("I19").Value = 965
$ws.Range("J19").Value = 1504.5
$ws.Range("K19").Value = 965
$ws.Range("L19").Value = 1504.5
$ws.Range("M19").Value = -795
$ws.Range("N19").Value = -1844.5
$ws.Range("H23").Value = 8741.75
$ws.Range("J23").Value = 4989
$ws.Range("L23").Value = 4989
$ws.Range("N23").Value = -5469
$ws.Range("H24").Value = 1273.2858
$ws.Range("I24").Value = 965
$ws.Range("J24").Value = 1504.5
$ws.Range("K24").Value = 965
$ws.Range("L24").Value = 1504.5
$ws.Range("M24").Value = -795
$ws.Range("N24").Value = -1844.5
$ws.Range("H25").Value = 1500
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H26").Value = 4257.5
$ws.Range("I26").Value = 4019
$ws.Range("J26").Value = 4337
$ws.Range("K26").Value = 4019
$ws.Range("L26").Value = 4337
$ws.Range("M26").Value = -3732
$ws.Range("N26").Value = -4911
$ws.Range("H27").Value = 8741.75
$ws.Range("J27").Value = 4989
$ws.Range("L27").Value = 4989
$ws.Range("N27").Value = -5373
$ws.Range("H31").Value = 165550.6
$ws.Range("I31").Value = 214843.81
$ws.Range("J31").Value = 43614.74
$ws.Range("K31").Value = 214843.81
$ws.Range("L31").Value = 43614.74
$ws.Range("M31").Value = -214548.81
$ws.Range("N31").Value = -44204.74
$ws.Range("H34").Value = 165550.6
$ws.Range("I34").Value = 214843.81
$ws.Range("J34").Value = 43614.74
$ws.Range("K34").Value = 214843.81
$ws.Range("L34").Value = 43614.74
$ws.Range("M34").Value = -214641.81
$ws.Range("N34").Value = -44018.74
$ws.Range("H35").Value = 809.375
$ws.Range("I35").Value = 1311.25
$ws.Range("K35").Value = 1311.25
$ws.Range("M35").Value = -1017.25
$ws.Range("H36").Value = 2997
$ws.Range("I36").Value = 2995
$ws.Range("K36").Value = 2995
$ws.Range("M36").Value = -2607
$ws.Range("H40").Value = 2997
$ws.Range("I40").Value = 2995
$ws.Range("K40").Value = 2995
$ws.Range("M40").Value = -2835
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H42").Value = 14494.5
$ws.Range("J42").Value = 14494.5
$ws.Range("L42").Value = 14494.5
$ws.Range("N42").Value = -15680.5
$ws.Range("H59").Value = 39857.145
$ws.Range("J59").Value = 39857.145
$ws.Range("L59").Value = 39857.145
$ws.Range("N59").Value = -42147.145
$ws.Range("H60").Value = 16425
$ws.Range("I60").Value = 12850
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 12850
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -12339
$ws.Range("N60").Value = -21022
$ws.Range("H107").Value = 1247.2195
$ws.Range("I107").Value = 625.8823
$ws.Range("J107").Value = 1687.3334
$ws.Range("K107").Value = 625.8823
$ws.Range("L107").Value = 1687.3334
$ws.Range("M107").Value = 1294.1177
$ws.Range("N107").Value = -5527.3334
$ws.Range("H113").Value = 1767.5
$ws.Range("I113").Value = 1981.3334
$ws.Range("J113").Value = 1639.2
$ws.Range("K113").Value = 1981.3334
$ws.Range("L113").Value = 1639.2
$ws.Range("M113").Value = 188.6666
$ws.Range("N113").Value = -5979.2
$ws.Range("H122").Value = 1320
$ws.Range("I122").Value = 1069.3334
$ws.Range("K122").Value = 3208.0002
$ws.Range("M122").Value = -758.0001999999999

# ---- Sheet: CUL (23 cell(s)) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 264277.88
$ws.Range("I4").Value = 972291.9
$ws.Range("J4").Value = 6818.227
$ws.Range("K4").Value = 2916875.7
$ws.Range("L4").Value = 20454.681
$ws.Range("M4").Value = -2916763.7
$ws.Range("N4").Value = -20678.681
$ws.Range("H16").Value = 7772.25
$ws.Range("J16").Value = 8962.666999999999
$ws.Range("L16").Value = 26888.001
$ws.Range("N16").Value = -27234.001
$ws.Range("H55").Value = 8209
$ws.Range("I55").Value = 6829.5713
$ws.Range("K55").Value = 20488.7139
$ws.Range("M55").Value = -20311.7139
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# ---- Sheet: GSM (19 cell(s)) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2500
$ws.Range("I5").Value = 4999
$ws.Range("K5").Value = 4999
$ws.Range("M5").Value = -4887
$ws.Range("H102").Value = 16069.25
$ws.Range("I102").Value = 3758.5
$ws.Range("K102").Value = 3758.5
$ws.Range("M102").Value = -2136.5
$ws.Range("H113").Value = 23326
$ws.Range("I113").Value = 23326
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 23326
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -21156
$ws.Range("N113").ClearContents()
$ws.Range("H123").Value = 39158.832
$ws.Range("J123").Value = 39158.832
$ws.Range("L123").Value = 39158.832
$ws.Range("N123").Value = -44058.832

# ---- Sheet: LTW (7 cell(s)) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 17000.5
$ws.Range("I2").Value = 22667
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 22667
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = -22555
$ws.Range("N2").Value = -225

# ---- Sheet: WVR (15 cell(s)) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 108207.75
$ws.Range("I2").Value = 431830
$ws.Range("J2").Value = 333.66666
$ws.Range("K2").Value = 431830
$ws.Range("L2").Value = 333.66666
$ws.Range("M2").Value = -431718
$ws.Range("N2").Value = -557.66666
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H136").Value = 271985.16
$ws.Range("I136").Value = 323762.88
$ws.Range("K136").Value = 971288.64
$ws.Range("M136").Value = -968738.64

